$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stok_jumlah values for existing rows
$ws.Range("B2").Value = 17
$ws.Range("B3").Value = 18

# Remove the last data row (row 4) entirely - shift cells up / delete row
$ws.Range("A4:E4").Delete()

$ws.Range("H7").Select()
